# Auto-generated edit script to update cryptos.xlsx data
# Applies the cell value changes described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text (General) representation
# by forcing Text number format before assigning values, so strings such as "1.00"
# or "0.600" are not silently converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "68.549.44"
$ws.Range("D3").Value = "2.694.34"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "599.09"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "160.40"
$ws.Range("E6").Value = "  +2.76%  "
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").Value = "2.693.99"
$ws.Range("E9").Value = "  +2.06%  "
$ws.Range("D10").Value = "0.141"
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("E13").Value = "  +2.61%  "
$ws.Range("D14").Value = "28.25"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").Value = "3.184.35"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "0.0000189"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "68.519.19"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "2.679.54"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").Value = "11.92"
$ws.Range("E19").Value = "  +5.24%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "366.96"
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "7.66"
$ws.Range("E21").Value = "  +3.76%  "
$ws.Range("E22").Value = "  +3.05%  "
$ws.Range("E23").Value = "  +2.37%  "
$ws.Range("D24").Value = "2.11"
$ws.Range("E24").Value = "  +2.80%  "
$ws.Range("D25").Value = "74.54"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "10.07"
$ws.Range("E27").Value = "  +3.97%  "
$ws.Range("D28").Value = "2.830.24"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").Value = "573.45"
$ws.Range("E31").Value = "  +3.52%  "
$ws.Range("E32").Value = "  +4.44%  "
$ws.Range("D33").Value = "8.23"
$ws.Range("E33").Value = "  +3.09%  "
$ws.Range("E34").Value = "  +5.71%  "
$ws.Range("E35").Value = "  +3.01%  "
$ws.Range("E36").Value = "  +6.60%  "
$ws.Range("D38").Value = "19.95"
$ws.Range("D39").Value = "160.61"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").Value = "0.379"
$ws.Range("E40").Value = "  +2.16%  "
$ws.Range("D41").Value = "1.92"
$ws.Range("E41").Value = "  +2.36%  "
$ws.Range("E42").Value = "  +2.05%  "
$ws.Range("D43").Value = "2.67"
$ws.Range("E43").Value = "  +2.30%  "
$ws.Range("E44").Value = "  +0.46%  "
$ws.Range("E46").Value = "  -6.22%  "
$ws.Range("D47").Value = "157.84"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").Value = "3.99"
$ws.Range("E48").Value = "  +7.23%  "
$ws.Range("E49").Value = "  +5.21%  "
$ws.Range("D50").Value = "0.600"
$ws.Range("E50").Value = "  +7.08%  "
$ws.Range("D51").Value = "22.06"
$ws.Range("E51").Value = "  +0.41%  "
